$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts old row 10 -> 11, old row 11 -> 12,
# and shifts their merged cells down automatically as well).
$ws.Rows("10:10").Insert()

# Seed the new row 10 with the same formatting / merges as the row-7 item
# template (copying a whole formatted+merged row is the most reliable way
# to reproduce the merge cells and per-cell styles for the new row).
$ws.Range("A7:Q7").Copy($ws.Range("A10:Q10"))

# Restore the explicit row heights the diff expects: the new item row keeps
# the 24.75 height, while the (shifted) totals row grows to 25.5.
$ws.Rows("10:10").RowHeight = 24.75
$ws.Rows("11:11").RowHeight = 25.5

# Fill in row 10 with the new item's data ("كالونا").
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "كالونا "
$ws.Range("H10").Value = "0:0"

# L10/P10 carry numeric-looking text ("0", "15.0000") but the template's
# number format for those columns is numeric, which would otherwise make
# Excel silently coerce the assigned text into a real number. Flip to a
# text format for the write, then restore the original numeric format so
# the cell keeps its original style index.
$fmtL10 = $ws.Range("L10").NumberFormat
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "0"
$ws.Range("L10").NumberFormat = $fmtL10

$ws.Range("N10").Value = "15.00"

$fmtP10 = $ws.Range("P10").NumberFormat
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "15.0000"
$ws.Range("P10").NumberFormat = $fmtP10

$ws.Range("Q10").Value = "1:0"

# Update the running total (P11, formerly P10) to include the new item.
$ws.Range("P11").Value = 202.48

# Update the generated timestamp footer (now A12) to the new save time.
$ws.Range("A12").Value = "Monday, 4 August, 2025 9:55 AM"
